$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(3408, 1).Value = 45565
$ws.Cells.Item(3408, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3408, 2).Value = 81111.17
$ws.Cells.Item(3408, 3).Value = 44.95
$ws.Cells.Item(3408, 4).Value = 5.58
$ws.Cells.Item(3408, 5).Value = 0.7

$ws.Cells.Item(3409, 1).Value = 45566
$ws.Cells.Item(3409, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3409, 2).Value = 81386.92999999999
$ws.Cells.Item(3409, 3).Value = 45.11
$ws.Cells.Item(3409, 4).Value = 5.6
$ws.Cells.Item(3409, 5).Value = 0.7

$ws.Cells.Item(3410, 1).Value = 45568
$ws.Cells.Item(3410, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3410, 2).Value = 79588.67
$ws.Cells.Item(3410, 3).Value = 44.11
$ws.Cells.Item(3410, 4).Value = 5.47
$ws.Cells.Item(3410, 5).Value = 0.72

$ws.Cells.Item(3411, 1).Value = 45569
$ws.Cells.Item(3411, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3411, 2).Value = 78846.7
$ws.Cells.Item(3411, 3).Value = 43.7
$ws.Cells.Item(3411, 4).Value = 5.42
$ws.Cells.Item(3411, 5).Value = 0.72

$ws.Cells.Item(3412, 1).Value = 45572
$ws.Cells.Item(3412, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3412, 2).Value = 77263.37
$ws.Cells.Item(3412, 3).Value = 42.82
$ws.Cells.Item(3412, 4).Value = 5.31
$ws.Cells.Item(3412, 5).Value = 0.74

$ws.Cells.Item(3413, 1).Value = 45573
$ws.Cells.Item(3413, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3413, 2).Value = 78929.56
$ws.Cells.Item(3413, 3).Value = 43.75
$ws.Cells.Item(3413, 4).Value = 5.43
$ws.Cells.Item(3413, 5).Value = 0.72

$ws.Cells.Item(3414, 1).Value = 45574
$ws.Cells.Item(3414, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3414, 2).Value = 79693.78
$ws.Cells.Item(3414, 3).Value = 44.17
$ws.Cells.Item(3414, 4).Value = 5.48
$ws.Cells.Item(3414, 5).Value = 0.71

$ws.Cells.Item(3415, 1).Value = 45575
$ws.Cells.Item(3415, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3415, 2).Value = 79468.87
$ws.Cells.Item(3415, 3).Value = 43.7
$ws.Cells.Item(3415, 4).Value = 5.39
$ws.Cells.Item(3415, 5).Value = 0.71

$ws.Cells.Item(3416, 1).Value = 45576
$ws.Cells.Item(3416, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3416, 2).Value = 79845.78999999999
$ws.Cells.Item(3416, 3).Value = 43.9
$ws.Cells.Item(3416, 4).Value = 5.42
$ws.Cells.Item(3416, 5).Value = 0.71

$ws.Cells.Item(3417, 1).Value = 45579
$ws.Cells.Item(3417, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3417, 2).Value = 80186.57000000001
$ws.Cells.Item(3417, 3).Value = 44.09
$ws.Cells.Item(3417, 4).Value = 5.44
$ws.Cells.Item(3417, 5).Value = 0.71

$ws.Cells.Item(3418, 1).Value = 45580
$ws.Cells.Item(3418, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3418, 2).Value = 80358.89
$ws.Cells.Item(3418, 3).Value = 44.15
$ws.Cells.Item(3418, 4).Value = 5.45
$ws.Cells.Item(3418, 5).Value = 0.71

$ws.Cells.Item(3419, 1).Value = 45581
$ws.Cells.Item(3419, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3419, 2).Value = 80168.23
$ws.Cells.Item(3419, 3).Value = 44.05
$ws.Cells.Item(3419, 4).Value = 5.44
$ws.Cells.Item(3419, 5).Value = 0.71

$ws.Cells.Item(3420, 1).Value = 45582
$ws.Cells.Item(3420, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3420, 2).Value = 78838.81
$ws.Cells.Item(3420, 3).Value = 43.24
$ws.Cells.Item(3420, 4).Value = 5.35
$ws.Cells.Item(3420, 5).Value = 0.72

$ws.Cells.Item(3421, 1).Value = 45583
$ws.Cells.Item(3421, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3421, 2).Value = 79085.86
$ws.Cells.Item(3421, 3).Value = 43.37
$ws.Cells.Item(3421, 4).Value = 5.37
$ws.Cells.Item(3421, 5).Value = 0.72

$ws.Cells.Item(3422, 1).Value = 45586
$ws.Cells.Item(3422, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3422, 2).Value = 77775.91
$ws.Cells.Item(3422, 3).Value = 42.73
$ws.Cells.Item(3422, 4).Value = 5.28
$ws.Cells.Item(3422, 5).Value = 0.73

$ws.Cells.Item(3423, 1).Value = 45587
$ws.Cells.Item(3423, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3423, 2).Value = 75748.27
$ws.Cells.Item(3423, 3).Value = 41.65
$ws.Cells.Item(3423, 4).Value = 5.14
$ws.Cells.Item(3423, 5).Value = 0.75

$ws.Cells.Item(3424, 1).Value = 45588
$ws.Cells.Item(3424, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3424, 2).Value = 76233.03999999999
$ws.Cells.Item(3424, 3).Value = 41.93
$ws.Cells.Item(3424, 4).Value = 5.17
$ws.Cells.Item(3424, 5).Value = 0.75

$ws.Cells.Item(3425, 1).Value = 45589
$ws.Cells.Item(3425, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3425, 2).Value = 75985.2
$ws.Cells.Item(3425, 3).Value = 41.09
$ws.Cells.Item(3425, 4).Value = 5.16
$ws.Cells.Item(3425, 5).Value = 0.75

$ws.Cells.Item(3426, 1).Value = 45590
$ws.Cells.Item(3426, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3426, 2).Value = 74539.92999999999
$ws.Cells.Item(3426, 3).Value = 40.11
$ws.Cells.Item(3426, 4).Value = 5.06
$ws.Cells.Item(3426, 5).Value = 0.76

$ws.Cells.Item(3427, 1).Value = 45593
$ws.Cells.Item(3427, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3427, 2).Value = 75158.39999999999
$ws.Cells.Item(3427, 3).Value = 40.31
$ws.Cells.Item(3427, 4).Value = 5.1
$ws.Cells.Item(3427, 5).Value = 0.76

$ws.Cells.Item(3428, 1).Value = 45594
$ws.Cells.Item(3428, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3428, 2).Value = 75852.47
$ws.Cells.Item(3428, 3).Value = 42.86
$ws.Cells.Item(3428, 4).Value = 5.15
$ws.Cells.Item(3428, 5).Value = 0.75

$ws.Cells.Item(3429, 1).Value = 45595
$ws.Cells.Item(3429, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3429, 2).Value = 75974.87
$ws.Cells.Item(3429, 3).Value = 42.92
$ws.Cells.Item(3429, 4).Value = 5.15
$ws.Cells.Item(3429, 5).Value = 0.75

$ws.Cells.Item(3430, 1).Value = 45596
$ws.Cells.Item(3430, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3430, 2).Value = 75670.5
$ws.Cells.Item(3430, 3).Value = 42.24
$ws.Cells.Item(3430, 4).Value = 5.13
$ws.Cells.Item(3430, 5).Value = 0.76

$ws.Cells.Item(3431, 1).Value = 45597
$ws.Cells.Item(3431, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3431, 2).Value = 76187.22
$ws.Cells.Item(3431, 3).Value = 42.53
$ws.Cells.Item(3431, 4).Value = 5.16
$ws.Cells.Item(3431, 5).Value = 0.75

$ws.Cells.Item(3432, 1).Value = 45600
$ws.Cells.Item(3432, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3432, 2).Value = 75239.67
$ws.Cells.Item(3432, 3).Value = 42.15
$ws.Cells.Item(3432, 4).Value = 5.1
$ws.Cells.Item(3432, 5).Value = 0.77

$ws.Cells.Item(3433, 1).Value = 45601
$ws.Cells.Item(3433, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3433, 2).Value = 75686
$ws.Cells.Item(3433, 3).Value = 42.4
$ws.Cells.Item(3433, 4).Value = 5.13
$ws.Cells.Item(3433, 5).Value = 0.77

$ws.Cells.Item(3434, 1).Value = 45602
$ws.Cells.Item(3434, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3434, 2).Value = 77358.92999999999
$ws.Cells.Item(3434, 3).Value = 43.38
$ws.Cells.Item(3434, 4).Value = 5.24
$ws.Cells.Item(3434, 5).Value = 0.74

$ws.Cells.Item(3435, 1).Value = 45603
$ws.Cells.Item(3435, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3435, 2).Value = 77026.2
$ws.Cells.Item(3435, 3).Value = 42.78
$ws.Cells.Item(3435, 4).Value = 5.22
$ws.Cells.Item(3435, 5).Value = 0.74

$ws.Cells.Item(3436, 1).Value = 45604
$ws.Cells.Item(3436, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3436, 2).Value = 76005.56
$ws.Cells.Item(3436, 3).Value = 42.2
$ws.Cells.Item(3436, 4).Value = 5.15
$ws.Cells.Item(3436, 5).Value = 0.77

$ws.Cells.Item(3437, 1).Value = 45607
$ws.Cells.Item(3437, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3437, 2).Value = 75333.53
$ws.Cells.Item(3437, 3).Value = 41.59
$ws.Cells.Item(3437, 4).Value = 5.1
$ws.Cells.Item(3437, 5).Value = 0.75

$ws.Cells.Item(3438, 1).Value = 45608
$ws.Cells.Item(3438, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3438, 2).Value = 74534.84
$ws.Cells.Item(3438, 3).Value = 41.13
$ws.Cells.Item(3438, 4).Value = 5.05
$ws.Cells.Item(3438, 5).Value = 0.77

$ws.Cells.Item(3439, 1).Value = 45609
$ws.Cells.Item(3439, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3439, 2).Value = 72569.97
$ws.Cells.Item(3439, 3).Value = 39.98
$ws.Cells.Item(3439, 4).Value = 4.92
$ws.Cells.Item(3439, 5).Value = 0.79

$ws.Cells.Item(3440, 1).Value = 45610
$ws.Cells.Item(3440, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3440, 2).Value = 72908.81
$ws.Cells.Item(3440, 3).Value = 40.14
$ws.Cells.Item(3440, 4).Value = 4.94
$ws.Cells.Item(3440, 5).Value = 0.8

$ws.Cells.Item(3441, 1).Value = 45614
$ws.Cells.Item(3441, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3441, 2).Value = 72911.08
$ws.Cells.Item(3441, 3).Value = 39.97
$ws.Cells.Item(3441, 4).Value = 4.94
$ws.Cells.Item(3441, 5).Value = 0.79

$ws.Cells.Item(3442, 1).Value = 45615
$ws.Cells.Item(3442, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3442, 2).Value = 73599.22
$ws.Cells.Item(3442, 3).Value = 39.97
$ws.Cells.Item(3442, 4).Value = 4.98
$ws.Cells.Item(3442, 5).Value = 0.8

$ws.Cells.Item(3443, 1).Value = 45617
$ws.Cells.Item(3443, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3443, 2).Value = 73379.44
$ws.Cells.Item(3443, 3).Value = 39.85
$ws.Cells.Item(3443, 4).Value = 4.97
$ws.Cells.Item(3443, 5).Value = 0.79

$ws.Cells.Item(3444, 1).Value = 45618
$ws.Cells.Item(3444, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3444, 2).Value = 74231.50999999999
$ws.Cells.Item(3444, 3).Value = 40.31
$ws.Cells.Item(3444, 4).Value = 5.03
$ws.Cells.Item(3444, 5).Value = 0.78

$ws.Cells.Item(3445, 1).Value = 45621
$ws.Cells.Item(3445, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3445, 2).Value = 75423.81
$ws.Cells.Item(3445, 3).Value = 40.96
$ws.Cells.Item(3445, 4).Value = 5.11
$ws.Cells.Item(3445, 5).Value = 0.77

$ws.Cells.Item(3446, 1).Value = 45622
$ws.Cells.Item(3446, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3446, 2).Value = 75442.5
$ws.Cells.Item(3446, 3).Value = 40.95
$ws.Cells.Item(3446, 4).Value = 5.11
$ws.Cells.Item(3446, 5).Value = 0.8100000000000001

$ws.Cells.Item(3447, 1).Value = 45623
$ws.Cells.Item(3447, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3447, 2).Value = 75925.48
$ws.Cells.Item(3447, 3).Value = 41.22
$ws.Cells.Item(3447, 4).Value = 5.08
$ws.Cells.Item(3447, 5).Value = 0.8100000000000001

$ws.Cells.Item(3448, 1).Value = 45624
$ws.Cells.Item(3448, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3448, 2).Value = 75963.84
$ws.Cells.Item(3448, 3).Value = 41.24
$ws.Cells.Item(3448, 4).Value = 5.08
$ws.Cells.Item(3448, 5).Value = 0.8100000000000001

$ws.Cells.Item(3449, 1).Value = 45625
$ws.Cells.Item(3449, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3449, 2).Value = 76087.8
$ws.Cells.Item(3449, 3).Value = 41.3
$ws.Cells.Item(3449, 4).Value = 5.09
$ws.Cells.Item(3449, 5).Value = 0.8100000000000001

$ws.Cells.Item(3450, 1).Value = 45628
$ws.Cells.Item(3450, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3450, 2).Value = 76908.39
$ws.Cells.Item(3450, 3).Value = 41.75
$ws.Cells.Item(3450, 4).Value = 5.14
$ws.Cells.Item(3450, 5).Value = 0.8

$ws.Cells.Item(3451, 1).Value = 45629
$ws.Cells.Item(3451, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3451, 2).Value = 77594.00999999999
$ws.Cells.Item(3451, 3).Value = 42.12
$ws.Cells.Item(3451, 4).Value = 5.19
$ws.Cells.Item(3451, 5).Value = 0.79

$ws.Cells.Item(3452, 1).Value = 45630
$ws.Cells.Item(3452, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3452, 2).Value = 78408.17999999999
$ws.Cells.Item(3452, 3).Value = 42.56
$ws.Cells.Item(3452, 4).Value = 5.24
$ws.Cells.Item(3452, 5).Value = 0.78

$ws.Cells.Item(3453, 1).Value = 45631
$ws.Cells.Item(3453, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3453, 2).Value = 78852.3
$ws.Cells.Item(3453, 3).Value = 42.8
$ws.Cells.Item(3453, 4).Value = 5.27
$ws.Cells.Item(3453, 5).Value = 0.78

$ws.Cells.Item(3454, 1).Value = 45632
$ws.Cells.Item(3454, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3454, 2).Value = 79207.2
$ws.Cells.Item(3454, 3).Value = 43
$ws.Cells.Item(3454, 4).Value = 5.29
$ws.Cells.Item(3454, 5).Value = 0.77

$ws.Cells.Item(3455, 1).Value = 45635
$ws.Cells.Item(3455, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3455, 2).Value = 79604.09
$ws.Cells.Item(3455, 3).Value = 43.21
$ws.Cells.Item(3455, 4).Value = 5.32
$ws.Cells.Item(3455, 5).Value = 0.77

$ws.Cells.Item(3456, 1).Value = 45636
$ws.Cells.Item(3456, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3456, 2).Value = 79788.47
$ws.Cells.Item(3456, 3).Value = 43.31
$ws.Cells.Item(3456, 4).Value = 5.33
$ws.Cells.Item(3456, 5).Value = 0.77

$ws.Cells.Item(3457, 1).Value = 45637
$ws.Cells.Item(3457, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3457, 2).Value = 80001.05
$ws.Cells.Item(3457, 3).Value = 43.43
$ws.Cells.Item(3457, 4).Value = 5.35
$ws.Cells.Item(3457, 5).Value = 0.77

$ws.Cells.Item(3458, 1).Value = 45638
$ws.Cells.Item(3458, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3458, 2).Value = 79635.08
$ws.Cells.Item(3458, 3).Value = 43.23
$ws.Cells.Item(3458, 4).Value = 5.32
$ws.Cells.Item(3458, 5).Value = 0.77

$ws.Cells.Item(3459, 1).Value = 45639
$ws.Cells.Item(3459, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3459, 2).Value = 79594.36
$ws.Cells.Item(3459, 3).Value = 43.21
$ws.Cells.Item(3459, 4).Value = 5.32
$ws.Cells.Item(3459, 5).Value = 0.77

$ws.Cells.Item(3460, 1).Value = 45642
$ws.Cells.Item(3460, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3460, 2).Value = 80203.53999999999
$ws.Cells.Item(3460, 3).Value = 43.54
$ws.Cells.Item(3460, 4).Value = 5.36
$ws.Cells.Item(3460, 5).Value = 0.76

$ws.Cells.Item(3461, 1).Value = 45643
$ws.Cells.Item(3461, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3461, 2).Value = 79743.23
$ws.Cells.Item(3461, 3).Value = 43.29
$ws.Cells.Item(3461, 4).Value = 5.33
$ws.Cells.Item(3461, 5).Value = 0.76

$ws.Cells.Item(3462, 1).Value = 45644
$ws.Cells.Item(3462, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3462, 2).Value = 79232.39
$ws.Cells.Item(3462, 3).Value = 43.01
$ws.Cells.Item(3462, 4).Value = 5.3
$ws.Cells.Item(3462, 5).Value = 0.77

$ws.Cells.Item(3463, 1).Value = 45645
$ws.Cells.Item(3463, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3463, 2).Value = 79007.06
$ws.Cells.Item(3463, 3).Value = 42.89
$ws.Cells.Item(3463, 4).Value = 5.28
$ws.Cells.Item(3463, 5).Value = 0.77

$ws.Cells.Item(3464, 1).Value = 45646
$ws.Cells.Item(3464, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3464, 2).Value = 76781.45
$ws.Cells.Item(3464, 3).Value = 41.68
$ws.Cells.Item(3464, 4).Value = 5.13
$ws.Cells.Item(3464, 5).Value = 0.79

$ws.Cells.Item(3465, 1).Value = 45649
$ws.Cells.Item(3465, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3465, 2).Value = 77032.64
$ws.Cells.Item(3465, 3).Value = 41.82
$ws.Cells.Item(3465, 4).Value = 5.15
$ws.Cells.Item(3465, 5).Value = 0.79

$ws.Cells.Item(3466, 1).Value = 45650
$ws.Cells.Item(3466, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3466, 2).Value = 76985.37
$ws.Cells.Item(3466, 3).Value = 41.79
$ws.Cells.Item(3466, 4).Value = 5.15
$ws.Cells.Item(3466, 5).Value = 0.79

$ws.Cells.Item(3467, 1).Value = 45652
$ws.Cells.Item(3467, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3467, 2).Value = 77076.87
$ws.Cells.Item(3467, 3).Value = 41.84
$ws.Cells.Item(3467, 4).Value = 5.15
$ws.Cells.Item(3467, 5).Value = 0.79

$ws.Cells.Item(3468, 1).Value = 45653
$ws.Cells.Item(3468, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3468, 2).Value = 76880.02
$ws.Cells.Item(3468, 3).Value = 41.73
$ws.Cells.Item(3468, 4).Value = 5.14
$ws.Cells.Item(3468, 5).Value = 0.79

$ws.Cells.Item(3469, 1).Value = 45656
$ws.Cells.Item(3469, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3469, 2).Value = 77163.27
$ws.Cells.Item(3469, 3).Value = 41.89
$ws.Cells.Item(3469, 4).Value = 5.16
$ws.Cells.Item(3469, 5).Value = 0.79

$ws.Cells.Item(3470, 1).Value = 45657
$ws.Cells.Item(3470, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3470, 2).Value = 77176.39
$ws.Cells.Item(3470, 3).Value = 42.9
$ws.Cells.Item(3470, 4).Value = 5.25
$ws.Cells.Item(3470, 5).Value = 0.78

$ws.Cells.Item(3471, 1).Value = 45658
$ws.Cells.Item(3471, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3471, 2).Value = 77515.66
$ws.Cells.Item(3471, 3).Value = 43.09
$ws.Cells.Item(3471, 4).Value = 5.28
$ws.Cells.Item(3471, 5).Value = 0.78

$ws.Cells.Item(3472, 1).Value = 45659
$ws.Cells.Item(3472, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3472, 2).Value = 78402.52
$ws.Cells.Item(3472, 3).Value = 43.58
$ws.Cells.Item(3472, 4).Value = 5.34
$ws.Cells.Item(3472, 5).Value = 0.77

$ws.Cells.Item(3473, 1).Value = 45660
$ws.Cells.Item(3473, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3473, 2).Value = 78163.48
$ws.Cells.Item(3473, 3).Value = 43.45
$ws.Cells.Item(3473, 4).Value = 5.32
$ws.Cells.Item(3473, 5).Value = 0.77

$ws.Cells.Item(3474, 1).Value = 45663
$ws.Cells.Item(3474, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3474, 2).Value = 76053.16
$ws.Cells.Item(3474, 3).Value = 42.28
$ws.Cells.Item(3474, 4).Value = 5.18
$ws.Cells.Item(3474, 5).Value = 0.79

$ws.Cells.Item(3475, 1).Value = 45664
$ws.Cells.Item(3475, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3475, 2).Value = 76730.91
$ws.Cells.Item(3475, 3).Value = 42.65
$ws.Cells.Item(3475, 4).Value = 5.22
$ws.Cells.Item(3475, 5).Value = 0.78

$ws.Cells.Item(3476, 1).Value = 45665
$ws.Cells.Item(3476, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3476, 2).Value = 75923.12
$ws.Cells.Item(3476, 3).Value = 42.2
$ws.Cells.Item(3476, 4).Value = 5.17
$ws.Cells.Item(3476, 5).Value = 0.79

$ws.Cells.Item(3477, 1).Value = 45666
$ws.Cells.Item(3477, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3477, 2).Value = 75215.17999999999
$ws.Cells.Item(3477, 3).Value = 41.81
$ws.Cells.Item(3477, 4).Value = 5.12
$ws.Cells.Item(3477, 5).Value = 0.8

$ws.Cells.Item(3478, 1).Value = 45667
$ws.Cells.Item(3478, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3478, 2).Value = 73649.82000000001
$ws.Cells.Item(3478, 3).Value = 40.94
$ws.Cells.Item(3478, 4).Value = 5.01
$ws.Cells.Item(3478, 5).Value = 0.82

$ws.Cells.Item(3479, 1).Value = 45670
$ws.Cells.Item(3479, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3479, 2).Value = 70687.74000000001
$ws.Cells.Item(3479, 3).Value = 39.3
$ws.Cells.Item(3479, 4).Value = 4.81
$ws.Cells.Item(3479, 5).Value = 0.85

$ws.Cells.Item(3480, 1).Value = 45671
$ws.Cells.Item(3480, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3480, 2).Value = 72423.02
$ws.Cells.Item(3480, 3).Value = 40.25
$ws.Cells.Item(3480, 4).Value = 4.93
$ws.Cells.Item(3480, 5).Value = 0.83

